# Applies the "update scripts wuth new tpm" edit to the Pomc-Mc5r LR-pairs sheet:
# expands the results table from 3 data rows to 8 (adds the MuSCs sending-cluster
# block and the Resolving-Mac target-cluster rows) and refreshes the TPM-derived
# metrics for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order: A..T =
# Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# Ligand-expressing cells, Ligand detection rate, Ligand average expression value,
# Ligand total expression value, Ligand derived specificity of average expression value,
# Ligand derived specificity of total expression value, Receptor-expressing cells,
# Receptor detection rate, Receptor average expression value, Receptor total expression value,
# Receptor derived specificity of average expression value, Receptor derived specificity of total expression value,
# Edge average expression weight, Edge total expression weight,
# Edge average expression derived specificity, Edge total expression derived specificity
$rows = @(
    @('ECs', 'Pomc', 'Mc5r', 'ECs', 3, 1, 0.5131323333333334, 1.539397, 0.9119558630037493, 0.9119558630037494, 3, 1, 1.011219333333333, 3.033658, 0.1971977600385235, 0.1971977600385236, 0.5188893360251111, 4.670004024226, 0.179835653438338, 0.179835653438338),
    @('ECs', 'Pomc', 'Mc5r', 'FAPs', 3, 1, 0.5131323333333334, 1.539397, 0.9119558630037493, 0.9119558630037494, 3, 1, 3.237665999999999, 9.712997999999999, 0.6313768555515021, 0.6313768555515022, 1.661351109134, 14.952159982206, 0.5757878251850637, 0.5757878251850639),
    @('ECs', 'Pomc', 'Mc5r', 'MuSCs', 3, 1, 0.5131323333333334, 1.539397, 0.9119558630037493, 0.9119558630037494, 3, 1, 0.7485606666666667, 2.245682, 0.1459767251808977, 0.1459767251808977, 0.3841106815282223, 3.456996133754, 0.1331243303908067, 0.1331243303908067),
    @('ECs', 'Pomc', 'Mc5r', 'Resolving-Mac', 3, 1, 0.5131323333333334, 1.539397, 0.9119558630037493, 0.9119558630037494, 1, 0.3333333333333333, 0.1304993333333333, 0.391498, 0.02544865922907654, 0.02544865922907655, 0.06696342741177778, 0.602670846706, 0.02320805398954083, 0.02320805398954083),
    @('MuSCs', 'Pomc', 'Mc5r', 'ECs', 1, 0.3333333333333333, 0.04954, 0.14862, 0.08804413699625062, 0.08804413699625063, 3, 1, 1.011219333333333, 3.033658, 0.1971977600385235, 0.1971977600385236, 0.05009580577333333, 0.45086225196, 0.01736210660018552, 0.01736210660018553),
    @('MuSCs', 'Pomc', 'Mc5r', 'FAPs', 1, 0.3333333333333333, 0.04954, 0.14862, 0.08804413699625062, 0.08804413699625063, 3, 1, 3.237665999999999, 9.712997999999999, 0.6313768555515021, 0.6313768555515022, 0.16039397364, 1.44354576276, 0.05558903036643839, 0.05558903036643841),
    @('MuSCs', 'Pomc', 'Mc5r', 'MuSCs', 1, 0.3333333333333333, 0.04954, 0.14862, 0.08804413699625062, 0.08804413699625063, 3, 1, 0.7485606666666667, 2.245682, 0.1459767251808977, 0.1459767251808977, 0.03708369542666667, 0.33375325884, 0.01285239479009098, 0.01285239479009098),
    @('MuSCs', 'Pomc', 'Mc5r', 'Resolving-Mac', 1, 0.3333333333333333, 0.04954, 0.14862, 0.08804413699625062, 0.08804413699625063, 1, 0.3333333333333333, 0.1304993333333333, 0.391498, 0.02544865922907654, 0.02544865922907655, 0.006464936973333333, 0.05818443276, 0.002240605239535713, 0.002240605239535713)
)

$rowCount = $rows.Count
$colCount = $rows[0].Count
$data = New-Object 'object[,]' $rowCount,$colCount
for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt $colCount; $c++) {
        $data[$r,$c] = $rows[$r][$c]
    }
}

$targetRange = $ws.Range("A2:T$(1 + $rowCount)")
$targetRange.Value = $data
